$d = $word.ActiveDocument
$d.Content.Find.Execute("https://doi.org/10.25573/serc.14714175", $false, $false, $false, $false, $false,
                         $true, 1, $false, "https://doi.org/10.25573/serc.14714175.v1", 2)
